# "Update v3 Code with Images"
#
# 1) Sheet1: a handful of cell-content updates
# 2) A new "Sheet2" (placed after Sheet1) holding the Insured/Claimant/
#    Payto-Receiver lookup list referenced by the new "Adjustor- Receiver
#    fraud Pair" sub-scenario comments
# 3) Restore the cursor/selection state that Excel persisted on save

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet1 content changes -------------------------------------------------

# Row 8 (merged B7:B9 "Multiple payments ... same cause" block): fix the
# capitalisation of the sub-scenario label.
$ws.Range("C8").Value = "Adjustor- Receiver fraud Pair"

# Row 9: status moved from In-Progress to Complete.
$ws.Range("E9").Value = "Complete"

# Rows 5 & 6 (merged B5:B6 block) Comments column: both now read "NA".
$ws.Range("I5").Value = "NA"
$ws.Range("I6").Value = "NA"

# --- New Sheet2: Insured / Claimant / Payto-Receiver reference list --------

$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Range("B2").Value = "Insured"
$ws2.Range("B3").Value = "Claimant"
$ws2.Range("B4").Value = "Payto/Receiver"
$null = $ws2.Range("B4").Select()

# --- Restore Sheet1 as the active sheet + its new selection ----------------

$null = $ws.Activate()
$null = $ws.Range("C7").Select()
